$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for columns B..O, rows 2..10 (row 1 is header, column A untouched)
$data = @{
    2  = @(20, 2, 18, 3, 9,  4, 13, 5, 8,  6, 22, 6, 10, 6)
    3  = @(26, 2, 19, 3, 18, 4, 6,  5, 7,  6, 17, 6, 17, 7)
    4  = @(26, 2, 20, 3, 15, 4, 11, 5, 5,  6, 19, 6, 21, 6)
    5  = @(19, 2, 13, 3, 7,  4, 5,  5, 6,  6, 25, 6, 23, 6)
    6  = @(18, 2, 5,  3, 7,  4, 16, 4, 13, 5, 6,  6, 15, 6)
    7  = @(16, 2, 36, 2, 5,  3, 17, 3, 11, 4, 30, 5, 17, 6)
    8  = @(26, 2, 23, 2, 5,  2, 19, 1, 5,  2, 28, 2, 3,  7)
    9  = @(14, 2, 18, 2, 9,  1, 5,  1, 4,  0, 35, 0, 15, 6)
    10 = @(18, 2, 17, 2, 2,  2, 4,  1, 3,  0, 1,  0, 3,  5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # column B is index 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
